$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2-4 are cyclically rotated: row2 <- row3, row3 <- row4, row4 <- row2
# (only columns A, B, E, F, G, H, Q, R actually change values; the rest stay identical)

# Capture original values for row 2 (the row whose data moves last, to row 4)
# NOTE: use .Value2 (not .Value) -- this runtime's .Value getter does not
# resolve to the underlying variant for Range in the same way .Value2 does.
$A2 = $ws.Range("A2").Value2
$B2 = $ws.Range("B2").Value2
$E2 = $ws.Range("E2").Value2
$F2 = $ws.Range("F2").Value2
$G2 = $ws.Range("G2").Value2
$H2 = $ws.Range("H2").Value2
$Q2 = $ws.Range("Q2").Value2
$R2 = $ws.Range("R2").Value2

# Capture original values for row 3 (moves up to row 2)
$A3 = $ws.Range("A3").Value2
$B3 = $ws.Range("B3").Value2
$E3 = $ws.Range("E3").Value2
$F3 = $ws.Range("F3").Value2
$G3 = $ws.Range("G3").Value2
$H3 = $ws.Range("H3").Value2
$Q3 = $ws.Range("Q3").Value2
$R3 = $ws.Range("R3").Value2

# Capture original values for row 4 (moves up to row 3)
$A4 = $ws.Range("A4").Value2
$B4 = $ws.Range("B4").Value2
$E4 = $ws.Range("E4").Value2
$F4 = $ws.Range("F4").Value2
$G4 = $ws.Range("G4").Value2
$H4 = $ws.Range("H4").Value2
$Q4 = $ws.Range("Q4").Value2
$R4 = $ws.Range("R4").Value2

# Write rotated values: row2 gets old row3 data
$ws.Range("A2").Value2 = $A3
$ws.Range("B2").Value2 = $B3
$ws.Range("E2").Value2 = $E3
$ws.Range("F2").Value2 = $F3
$ws.Range("G2").Value2 = $G3
$ws.Range("H2").Value2 = $H3
$ws.Range("Q2").Value2 = $Q3
$ws.Range("R2").Value2 = $R3

# row3 gets old row4 data
$ws.Range("A3").Value2 = $A4
$ws.Range("B3").Value2 = $B4
$ws.Range("E3").Value2 = $E4
$ws.Range("F3").Value2 = $F4
$ws.Range("G3").Value2 = $G4
$ws.Range("H3").Value2 = $H4
$ws.Range("Q3").Value2 = $Q4
$ws.Range("R3").Value2 = $R4

# row4 gets old row2 data
$ws.Range("A4").Value2 = $A2
$ws.Range("B4").Value2 = $B2
$ws.Range("E4").Value2 = $E2
$ws.Range("F4").Value2 = $F2
$ws.Range("G4").Value2 = $G2
$ws.Range("H4").Value2 = $H2
$ws.Range("Q4").Value2 = $Q2
$ws.Range("R4").Value2 = $R2
